$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(3923, 4028, 4312, 4312, 4366, 4535, 4535, 4535, 4535, 4535, 4535, 4592, 4592, 4592)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $values[$i]
}
